$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mistral")

$data = @(
    @('preds_ns5_ws200_mc1000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_lenNone_gblFalse', 2624.061439999999),
    @('preds_ns5_ws200_mc750_ea1.002_snks0_hopf_True_type_max_fused_opt_qcache_new_burst_lenNone_gblFalse', 98.43507200000001),
    @('preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_prof_lenNone_gblFalse', 1343.823872),
    @('preds_ns5_ws200_mc4000_ea1.002_snks0_hopf_True_type_max_fused_opt_qcache_new_burst_lenNone_gblFalse', 10488.38144),
    @('preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_opt_lenNone_gblFalse', 42281.28768),
    @('preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_snapkv_prof_lenNone_gblFalse', 5255.987199999999),
    @('preds_ns5_ws200_mc2000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_lenNone_gblFalse', 10491.00288000001),
    @('preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_h2o_prof_qcache_lenNone_gblFalse', 396.816384),
    @('preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_lenNone_gblFalse', 10488.38144),
    @('preds_ns1_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_prof_qcache_lenNone_gblFalse', 1344.274432),
    @('preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_snapkv_opt_lenNone_gblFalse', 84484.292608),
    @('preds_ns10_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_rerun_lenNone_gblFalse', 53752.95488000003),
    @('preds_ns5_ws200_mc2000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_sl3_lenNone_gblFalse', 21286.223872),
    @('preds_ns5_ws200_mc1000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_new_lenNone_gblFalse', 2624.061439999999),
    @('preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_opt_qcache_lenNone_gblFalse', 42281.28768),
    @('preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_opt_lenNone_gblFalse', 28220.30131200001),
    @('preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_snapkv_opt_qcache_lenNone_gblFalse', 84484.292608),
    @('preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_snapkv_prof_qcache_lenNone_gblFalse', 1763.704832),
    @('preds_ns5_ws200_mc1000_ea1.0_snks0_hopf_False_type_max_fused_lenNone_gblFalse', 21121.073152),
    @('preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_max_fused_prof_lenNone_gblFalse', 1261.19936),
    @('preds_ns10_ws32_mc4000_ea1.0_snks0_hopf_True_type_max_fused_rerun_lenNone_gblFalse', 26220.95360000002),
    @('preds_ns10_ws200_mc1000_ea1.0_snks0_hopf_True_type_snapkv_rerun_lenNone_gblFalse', 169679.519744),
    @('preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_opt_qcache_new_lenNone_gblFalse', 42281.28768),
    @('preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_snapkv_opt_qcache_new_lenNone_gblFalse', 84484.292608),
    @('preds_ns1_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_prof_lenNone_gblFalse', 0),
    @('preds_ns5_ws200_mc750_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_new_burst_lenNone_gblFalse', 98.43507200000001),
    @('preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_max_fused_prof_qcache_lenNone_gblFalse', 492.17536),
    @('preds_ns10_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_rerun_lenNone_gblFalse', 84562.57535999996),
    @('preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_new_lenNone_gblFalse', 10488.38144),
    @('preds_ns10_ws200_mc1000_ea1.0_snks0_hopf_True_type_h2o_rerun_lenNone_gblFalse', 21156.49535999999),
    @('preds_ns10_ws200_mc4000_ea1.0_snks0_hopf_True_type_sum_fused_rerun_lenNone_gblFalse', 53752.95488000003)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}

# Make sure the new label cells (A25:A32) pick up the same formatting
# (bold, centered, bordered) already used by the existing label column.
$fmtSource = $ws.Range("A2")
$fmtSource.Copy()
$ws.Range("A25:A32").PasteSpecial(-4122)
$excel.CutCopyMode = 0
